# Updates cryptos list data (prices and 1h volume changes) per commit
# "Updated cryptos list on Sat Apr 15 12:36:15 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the affected ranges to Text format first so that numeric-looking
# strings (e.g. "53.80", "1.002", "0.00001130") are stored verbatim as text
# rather than being coerced into floating point numbers.
$ws.Range("B40:E41").NumberFormat = "@"
$ws.Range("B49:E51").NumberFormat = "@"
$ws.Range("D2:E39").NumberFormat = "@"
$ws.Range("D42:E48").NumberFormat = "@"

$ws.Range("D2").Value = '30.477.54'
$ws.Range("E2").Value = '  -1.11%  '
$ws.Range("D3").Value = '2.108.25'
$ws.Range("E3").Value = '  -0.41%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '333.88'
$ws.Range("E5").Value = '  -0.15%  '
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("D7").Value = '0.5258'
$ws.Range("E7").Value = '  -1.29%  '
$ws.Range("D8").Value = '0.4583'
$ws.Range("E8").Value = '  +4.62%  '
$ws.Range("D9").Value = '53.80'
$ws.Range("E9").Value = '  +14.47%  '
$ws.Range("D10").Value = '0.09001'
$ws.Range("E10").Value = '  +0.11%  '
$ws.Range("D11").Value = '1.182'
$ws.Range("E11").Value = '  +0.44%  '
$ws.Range("D12").Value = '24.46'
$ws.Range("E12").Value = '  -2.28%  '
$ws.Range("D13").Value = '2.094.34'
$ws.Range("E13").Value = '  -0.96%  '
$ws.Range("D14").Value = '6.789'
$ws.Range("E14").Value = '  +0.46%  '
$ws.Range("D15").Value = '7.841'
$ws.Range("E15").Value = '  +0.43%  '
$ws.Range("D16").Value = '96.83'
$ws.Range("E16").Value = '  -0.22%  '
$ws.Range("D17").Value = '1.003'
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("D18").Value = '0.00001130'
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("E19").Value = '  -0.90%  '
$ws.Range("D20").Value = '19.48'
$ws.Range("E20").Value = '  +1.86%  '
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("D22").Value = '6.322'
$ws.Range("E22").Value = '  -0.35%  '
$ws.Range("D23").Value = '30.522.23'
$ws.Range("E23").Value = '  -1.22%  '
$ws.Range("E24").Value = '  +0.67%  '
$ws.Range("D25").Value = '2.347'
$ws.Range("E25").Value = '  +2.38%  '
$ws.Range("D26").Value = '2.343.98'
$ws.Range("E26").Value = '  -0.85%  '
$ws.Range("D27").Value = '22.38'
$ws.Range("E27").Value = '  -1.69%  '
$ws.Range("D28").Value = '2.585'
$ws.Range("E28").Value = '  -0.58%  '
$ws.Range("D29").Value = '163.61'
$ws.Range("E29").Value = '  +0.16%  '
$ws.Range("D30").Value = '133.03'
$ws.Range("E30").Value = '  -0.40%  '
$ws.Range("D31").Value = '1.198'
$ws.Range("E31").Value = '  +0.57%  '
$ws.Range("E32").Value = '  -1.06%  '
$ws.Range("D33").Value = '1.676'
$ws.Range("E33").Value = '  +7.58%  '
$ws.Range("D34").Value = '6.149'
$ws.Range("E34").Value = '  -1.02%  '
$ws.Range("D35").Value = '3.927'
$ws.Range("E35").Value = '  -3.21%  '
$ws.Range("D36").Value = '10.45'
$ws.Range("E36").Value = '  +8.27%  '
$ws.Range("D37").Value = '0.02577'
$ws.Range("E37").Value = '  -1.06%  '
$ws.Range("D38").Value = '0.06837'
$ws.Range("E38").Value = '  +0.68%  '
$ws.Range("D39").Value = '5.575'
$ws.Range("E39").Value = '  +0.45%  '
$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").Value = '12.77'
$ws.Range("E40").Value = '  -0.17%  '
$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").Value = '0.2292'
$ws.Range("E41").Value = '  -0.37%  '
$ws.Range("D42").Value = '0.6909'
$ws.Range("E42").Value = '  +0.74%  '
$ws.Range("D43").Value = '1.243'
$ws.Range("E43").Value = '  -0.27%  '
$ws.Range("D44").Value = '2.346'
$ws.Range("E45").Value = '  +0.10%  '
$ws.Range("D46").Value = '14.07'
$ws.Range("E46").Value = '  -0.24%  '
$ws.Range("D47").Value = '0.6392'
$ws.Range("E47").Value = '  -0.89%  '
$ws.Range("D48").Value = '3.655'
$ws.Range("E48").Value = '  -0.40%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.00000000353'
$ws.Range("E49").Value = '  +27.10%  '
$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D50").Value = '1.244'
$ws.Range("E50").Value = '  -1.80%  '
$ws.Range("B51").Value = 'WEMIXTOKEN'
$ws.Range("C51").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D51").Value = '1.218'
$ws.Range("E51").Value = '  +0.66%  '

# Restore default (Normal) style so no stray number-format styling is
# left attached to the cells - matches the original (unstyled) cells.
$ws.Range("B40:E41").Style = "Normal"
$ws.Range("B49:E51").Style = "Normal"
$ws.Range("D2:E39").Style = "Normal"
$ws.Range("D42:E48").Style = "Normal"

